$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "ngClass e ngStyle" paragraph (currently paragraph #18, 1-based):
#    - make the "ngClass e ngStyle:" label bold
#    - make the trailing "." bold
#    - remove the single-underline paragraph-mark formatting
# ---------------------------------------------------------------------------
$pClassStyle = $d.Paragraphs.Item(18)
$pClassStyle.Range.Font.Underline = 0

$labelStart = $pClassStyle.Range.Start
$labelEnd = $labelStart + [int]("ngClass e ngStyle:".Length)
$labelRange = $d.Range($labelStart, $labelEnd)
$labelRange.Font.Bold = 1

$paraEnd = $pClassStyle.Range.End
$dotStart = $paraEnd - 2
$dotRange = $d.Range($dotStart, $paraEnd - 1)
$dotRange.Font.Bold = 1

# ---------------------------------------------------------------------------
# 2) Turn the following empty paragraph (#19) into the new
#    "ngShow, ngHide ... ngIf ..." paragraph, bold where appropriate, and
#    put the _GoBack bookmark back in (after "(tags pequenas").
# ---------------------------------------------------------------------------
$pShow = $d.Paragraphs.Item(19)
$pShow.Range.InsertBefore("ngShow, ngHide (tags pequenas) e ngIf (melhor para performance - associar ngsource): Exibindo um elemento condicionalmente. ngIf interage com a DOM já os outros não.")

$showStart = $pShow.Range.Start

# bold: "ngShow, ngHide (tags pequenas) e ngIf "
$b1Len = [int]("ngShow, ngHide (tags pequenas) e ngIf ".Length)
$b1 = $d.Range($showStart, $showStart + $b1Len)
$b1.Font.Bold = 1

# plain: "(melhor para performance - associar ngsource)"
$plain1Len = [int]("(melhor para performance - associar ngsource)".Length)
$plain1Start = $showStart + $b1Len
$plain1 = $d.Range($plain1Start, $plain1Start + $plain1Len)
$plain1.Font.Bold = 0

# bold: ":"
$boldColonStart = $plain1Start + $plain1Len
$boldColon = $d.Range($boldColonStart, $boldColonStart + 1)
$boldColon.Font.Bold = 1

# plain: " Exibindo um elemento condicionalmente. ngIf interage com a DOM já os outros não."
$restStart = $boldColonStart + 1
$restEnd = $pShow.Range.End - 1
$restRange = $d.Range($restStart, $restEnd)
$restRange.Font.Bold = 0

# Set the paragraph mark itself to bold (matches the pPr/rPr of the target)
$pShow.Range.Font.Underline = 0

# place the _GoBack bookmark right after "(tags pequenas"
$bmPos = $showStart + [int]("ngShow, ngHide (tags pequenas".Length)
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# 3) Insert a new blank paragraph (with the single-underline paragraph mark)
#    right after the new "ngShow..." paragraph, to keep the original spacing
#    that used to follow "ngClass e ngStyle".
# ---------------------------------------------------------------------------
$pShow2 = $d.Paragraphs.Item(19)
$pShow2.Range.InsertParagraphAfter()

$pNewBlank = $d.Paragraphs.Item(20)
$pNewBlank.Range.Font.Bold = 0
$pNewBlank.Range.Font.Underline = 1
